$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new column before column N ---
# (Loan RBI / Variable Instalments change: a new blank column is inserted
#  between "In Advance" (M) and "Late" (old N), pushing "Late", the blank
#  "heading" separator column, and "Outstanding" one column to the right.)
$ws = $wb.Worksheets.Item("Repayment schedule")

# Capture the width of the column immediately to the left (M) so the newly
# inserted column N can inherit the same width, just like Excel does when
# inserting a column.
$mWidth = $ws.Columns("M").ColumnWidth

$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# --- Make "Repayment schedule" the active sheet/tab, with P7 selected ---
$ws.Activate()
$ws.Range("P7").Select() | Out-Null
